$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header text updates (rich-text shared strings) - edit only the
#    sub-string run that changed, leaving the surrounding runs untouched.
# ---------------------------------------------------------------------------

# A8: "Volume 30   Number  25" -> "...26"
$a8 = $ws.Range("A8")
$a8Text = $a8.Text
$a8.Characters($a8Text.Length - 1, 2).Text = "26"

# C9: "Report Covering the Week  6/19/2023  Through  6/25/2023"
#     -> "...6/26/2023  Through  7/2/2023"
$c9 = $ws.Range("C9")
$c9.Characters(27, 9).Text = "6/26/2023"
$c9.Characters(47, 9).Text = "7/2/2023"

# ---------------------------------------------------------------------------
# 2. Helper to flip a currently-numeric cell into a text/shared-string cell
#    (used where the new value is "0" or "***.*"), while keeping the same
#    visual style as a same-styled donor cell.
# ---------------------------------------------------------------------------
function Set-TextCell($cellRef, $text, $donorRef) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $ws.Range($donorRef).Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# 3. Helper to flip a currently-text cell into a numeric cell with a given
#    number format string (matches the target style's numFmt).
# ---------------------------------------------------------------------------
function Set-NumberCell($cellRef, $value, $numberFormat) {
    $cell = $ws.Range($cellRef)
    $cell.Value = $value
    $cell.NumberFormat = $numberFormat
}

# ---------------------------------------------------------------------------
# Row 15: D15/E15 go from "N/A" / "***.*" text to real numbers
# ---------------------------------------------------------------------------
Set-NumberCell "D15" 1 "#,##0"
Set-NumberCell "E15" -100 '#,##0.0;"-"#,##0.0'

# ---------------------------------------------------------------------------
# Row 26: D26/E26 go from "N/A" / "***.*" text to real numbers
# ---------------------------------------------------------------------------
Set-NumberCell "D26" 1 "#,##0"
Set-NumberCell "E26" -100 '#,##0.0;"-"#,##0.0'

# ---------------------------------------------------------------------------
# Row 27: C27 goes from number to "N/A" text
# ---------------------------------------------------------------------------
Set-TextCell "C27" "0" "D14"

# ---------------------------------------------------------------------------
# Row 30: G30/H30 go from numbers to "N/A" / "***.*" text
# ---------------------------------------------------------------------------
Set-TextCell "G30" "0" "F30"
Set-TextCell "H30" "***.*" "N22"

# ---------------------------------------------------------------------------
# 4. Plain numeric value updates (style/type unchanged)
# ---------------------------------------------------------------------------

# Row 14
$ws.Range("N14").Value = -68.421052631578

# Row 15 (remaining numeric-only changes)
$ws.Range("J15").Value = 11
$ws.Range("K15").Value = -18.181818181818
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = -73.529411764705

# Row 16
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -28.571428571428
$ws.Range("F16").Value = 19
$ws.Range("G16").Value = 23
$ws.Range("H16").Value = -17.391304347826
$ws.Range("I16").Value = 91
$ws.Range("J16").Value = 98
$ws.Range("K16").Value = -7.142857142857
$ws.Range("L16").Value = 5.813953488372
$ws.Range("M16").Value = -19.469026548672
$ws.Range("N16").Value = -78.281622911694

# Row 17
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 16.666666666666
$ws.Range("F17").Value = 38
$ws.Range("G17").Value = 45
$ws.Range("H17").Value = -15.555555555555
$ws.Range("I17").Value = 190
$ws.Range("J17").Value = 184
$ws.Range("K17").Value = 3.260869565217
$ws.Range("L17").Value = 2.702702702702
$ws.Range("M17").Value = 95.876288659793
$ws.Range("N17").Value = -35.374149659863

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -40
$ws.Range("I18").Value = 58
$ws.Range("J18").Value = 85
$ws.Range("K18").Value = -31.764705882352
$ws.Range("L18").Value = 18.367346938775
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = -71.428571428571

# Row 19
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 350
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 208
$ws.Range("J19").Value = 225
$ws.Range("K19").Value = -7.555555555555
$ws.Range("L19").Value = 22.352941176470
$ws.Range("M19").Value = 67.741935483871
$ws.Range("N19").Value = -31.353135313531

# Row 20
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 50
$ws.Range("L20").Value = 50
$ws.Range("M20").Value = 90.909090909090
$ws.Range("N20").Value = -77.419354838709

# Row 21
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 35.294117647058
$ws.Range("F21").Value = 104
$ws.Range("G21").Value = 118
$ws.Range("H21").Value = -11.864406779661
$ws.Range("I21").Value = 604
$ws.Range("J21").Value = 648
$ws.Range("K21").Value = -6.790123456790
$ws.Range("L21").Value = 13.747645951035
$ws.Range("M21").Value = 52.911392405063
$ws.Range("N21").Value = -58.573388203017

# Row 22
$ws.Range("L22").Value = -33.333333333333
$ws.Range("M22").Value = -33.333333333333

# Row 23
$ws.Range("D23").Value = 8
$ws.Range("E23").Value = -12.5
$ws.Range("F23").Value = 29
$ws.Range("G23").Value = 36
$ws.Range("H23").Value = -19.444444444444
$ws.Range("I23").Value = 188
$ws.Range("J23").Value = 193
$ws.Range("K23").Value = -2.590673575129
$ws.Range("L23").Value = -6.930693069306
$ws.Range("M23").Value = 70.909090909090

# Row 24
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = -27.272727272727
$ws.Range("I24").Value = 442
$ws.Range("J24").Value = 404
$ws.Range("K24").Value = 9.405940594059
$ws.Range("L24").Value = 16.931216931216
$ws.Range("M24").Value = 41.666666666666

# Row 25
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 61
$ws.Range("G25").Value = 45
$ws.Range("H25").Value = 35.555555555555
$ws.Range("I25").Value = 300
$ws.Range("J25").Value = 297
$ws.Range("K25").Value = 1.010101010101
$ws.Range("L25").Value = 12.781954887218
$ws.Range("M25").Value = -17.127071823204

# Row 26 (remaining numeric-only changes)
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 1
$ws.Range("J26").Value = 14
$ws.Range("K26").Value = 7.142857142857
$ws.Range("L26").Value = 7.142857142857

# Row 27 (remaining numeric-only changes)
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = -16.666666666666
$ws.Range("J27").Value = 33
$ws.Range("K27").Value = -12.121212121212

# Row 28
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 100
$ws.Range("L28").Value = -44.444444444444
$ws.Range("M28").Value = -50
$ws.Range("N28").Value = -77.777777777777

# Row 29
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 100
$ws.Range("L29").Value = -47.058823529411
$ws.Range("M29").Value = -47.058823529411
$ws.Range("N29").Value = -78.571428571428

Write-Host "Edit complete"
